$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "42.822.74"
$ws.Cells.Item(2, 5).Value = "  -0.28%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.249.69"
$ws.Cells.Item(3, 5).Value = "  +0.47%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.01"
$ws.Cells.Item(4, 5).Value = "  +0.28%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "113.45"
$ws.Cells.Item(5, 5).Value = "  -0.51%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "295.77"
$ws.Cells.Item(6, 5).Value = "  +7.02%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.03%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.32%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.606"
$ws.Cells.Item(9, 5).Value = "  -0.44%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "44.09"
$ws.Cells.Item(10, 5).Value = "  -4.81%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0924"
$ws.Cells.Item(11, 5).Value = "  -0.69%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "54.30"
$ws.Cells.Item(12, 5).Value = "  +0.06%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "8.94"
$ws.Cells.Item(13, 5).Value = "  -1.27%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "1.06"
$ws.Cells.Item(14, 5).Value = "  +21.43%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -1.21%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "15.10"
$ws.Cells.Item(16, 5).Value = "  -1.39%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.586.54"
$ws.Cells.Item(17, 5).Value = "  +0.36%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.246.14"
$ws.Cells.Item(18, 5).Value = "  -0.09%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "42.695.69"
$ws.Cells.Item(19, 5).Value = "  -0.70%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.23"
$ws.Cells.Item(20, 5).Value = "  +6.81%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -0.90%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "74.60"
$ws.Cells.Item(22, 5).Value = "  +3.43%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.51"
$ws.Cells.Item(23, 5).Value = "  +16.73%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.41"
$ws.Cells.Item(24, 5).Value = "  +2.93%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "249.25"
$ws.Cells.Item(25, 5).Value = "  +7.50%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "8.97"
$ws.Cells.Item(26, 5).Value = "  -3.71%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.91%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "11.55"
$ws.Cells.Item(28, 5).Value = "  -5.20%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -1.14%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "175.83"
$ws.Cells.Item(30, 5).Value = "  +1.47%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "37.72"
$ws.Cells.Item(31, 5).Value = "  -6.72%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "22.01"
$ws.Cells.Item(32, 5).Value = "  +4.38%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -4.03%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0887"
$ws.Cells.Item(34, 5).Value = "  -0.90%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.70"
$ws.Cells.Item(35, 5).Value = "  +2.36%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.07"
$ws.Cells.Item(36, 5).Value = "  +9.21%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.26"
$ws.Cells.Item(37, 5).Value = "  -2.55%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.127"
$ws.Cells.Item(38, 5).Value = "  -0.07%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.52%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -1.43%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.42"
$ws.Cells.Item(41, 5).Value = "  -5.50%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "72.13"
$ws.Cells.Item(42, 5).Value = "  +1.31%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.232"
$ws.Cells.Item(43, 5).Value = "  -0.59%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.15%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "12.50"
$ws.Cells.Item(45, 5).Value = "  -5.05%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.67%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "5.50"
$ws.Cells.Item(47, 5).Value = "  -2.90%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.23%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "105.35"
$ws.Cells.Item(49, 5).Value = "  +4.92%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.62"
$ws.Cells.Item(50, 5).Value = "  +2.39%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "71.54"
$ws.Cells.Item(51, 5).Value = "  +0.83%  "
